# "updates methods and rewrites metadata"
#
# This reproduces the content edits made to the Feather River snorkel
# metadata workbook: the dataset title/abstract was rewritten, the
# keyword list was refreshed (old run-timing / life-stage keywords swapped
# for species + habitat + place keywords), and a second taxon (steelhead)
# was added to the taxonomic coverage sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) title sheet -- rewrite the dataset title (short_name is unchanged)
# ------------------------------------------------------------------
$titleSheet = $wb.Worksheets.Item("title")
$titleSheet.Range("A2").Value = "Distribution and habitat use of juvenile Feather River salmonids"

# ------------------------------------------------------------------
# 2) keyword_set sheet -- replace the old life-stage/run keywords with
#    the new species / habitat / place keywords used by the rewritten
#    metadata. Taxonomic + place keywords (Oncorhynchus tshawytscha,
#    California, Central Valley) stay, shifted up, and "Feather River"
#    is appended as a new keyword.
# ------------------------------------------------------------------
$keywordSheet = $wb.Worksheets.Item("keyword_set")
$keywordSheet.Range("A2").Value = "chinook"
$keywordSheet.Range("A3").Value = "habitat"
$keywordSheet.Range("A4").Value = "Oncorhynchus tshawytscha"
$keywordSheet.Range("A5").Value = "California"
$keywordSheet.Range("A6").Value = "Central Valley"
$keywordSheet.Range("A7").Value = "Feather River"

# ------------------------------------------------------------------
# 3) taxonomic_coverage sheet -- add steelhead as a second covered taxon
# ------------------------------------------------------------------
$taxSheet = $wb.Worksheets.Item("taxonomic_coverage")
$taxSheet.Range("A3").Value = "steelhead"
$taxSheet.Range("A3").Select()

# ------------------------------------------------------------------
# 4) taxonomic_coverage becomes the active/last-viewed sheet
# ------------------------------------------------------------------
$taxSheet.Activate()
